$wb = $excel.ActiveWorkbook

# --- "About" sheet: update the "last updated" date in C1 (2024-03-15 -> 2024-03-28) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- "RAF-capacity" sheet: bump the RAF values for the two hydrogen technologies to 1 ---
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# Widen column A slightly on the RAF-capacity sheet
$wsCapacity.Columns("A:A").ColumnWidth = 29.04296875

# --- Window / view state: RAF-capacity becomes the active tab, scrolled & zoomed ---
$wsCapacity.Activate()
$wsCapacity.Range("B25").Select()
$excel.ActiveWindow.Zoom = 80
$excel.ActiveWindow.ScrollRow = 14
